$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.008646413965192
$ws.Cells.Item(2, 4).Value = 1.011292129856397
$ws.Cells.Item(2, 5).Value = 1.011191896199557
$ws.Cells.Item(2, 6).Value = 1.01229742457813
$ws.Cells.Item(2, 9).Value = 1.022827583285258
$ws.Cells.Item(2, 10).Value = 1.013911034312758
$ws.Cells.Item(2, 11).Value = 1.01415974198265
$ws.Cells.Item(2, 12).Value = 1.01405981078793
$ws.Cells.Item(2, 13).Value = 1.015162006767323
$ws.Cells.Item(2, 14).Value = 1.008810107016128

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.010021944816173
$ws.Cells.Item(3, 4).Value = 1.012578359252706
$ws.Cells.Item(3, 5).Value = 1.012366064095743
$ws.Cells.Item(3, 6).Value = 1.014951094925411
$ws.Cells.Item(3, 9).Value = 1.023108887735401
$ws.Cells.Item(3, 10).Value = 1.0149163029505
$ws.Cells.Item(3, 11).Value = 1.015248885006025
$ws.Cells.Item(3, 12).Value = 1.0150371834
$ws.Cells.Item(3, 13).Value = 1.017615005122221
$ws.Cells.Item(3, 14).Value = 1.009140587697577

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.010903280765557
$ws.Cells.Item(4, 4).Value = 1.013402278278667
$ws.Cells.Item(4, 5).Value = 1.013118706724255
$ws.Cells.Item(4, 6).Value = 1.016618885653574
$ws.Cells.Item(4, 9).Value = 1.023277673179666
$ws.Cells.Item(4, 10).Value = 1.015557744074275
$ws.Cells.Item(4, 11).Value = 1.015944766446849
$ws.Cells.Item(4, 12).Value = 1.015661947204082
$ws.Cells.Item(4, 13).Value = 1.019152871948848
$ws.Cells.Item(4, 14).Value = 1.009351445949817

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.011271741234859
$ws.Cells.Item(5, 4).Value = 1.013746686941099
$ws.Cells.Item(5, 5).Value = 1.013433442011817
$ws.Cells.Item(5, 6).Value = 1.017308393167433
$ws.Cells.Item(5, 9).Value = 1.023345486493278
$ws.Cells.Item(5, 10).Value = 1.015825274874751
$ws.Cells.Item(5, 11).Value = 1.016235224791662
$ws.Cells.Item(5, 12).Value = 1.015922792120972
$ws.Cells.Item(5, 13).Value = 1.019787733302648
$ws.Cells.Item(5, 14).Value = 1.009439386693721

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.011333487961706
$ws.Cells.Item(6, 4).Value = 1.013804400276096
$ws.Cells.Item(6, 5).Value = 1.013486190010977
$ws.Cells.Item(6, 6).Value = 1.017423487248807
$ws.Cells.Item(6, 9).Value = 1.023356688957821
$ws.Cells.Item(6, 10).Value = 1.015870070396353
$ws.Cells.Item(6, 11).Value = 1.016283872361193
$ws.Cells.Item(6, 12).Value = 1.015966484029825
$ws.Cells.Item(6, 13).Value = 1.019893650290163
$ws.Cells.Item(6, 14).Value = 1.009454111337598

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.01090821217382
$ws.Cells.Item(7, 4).Value = 1.013406887968694
$ws.Cells.Item(7, 5).Value = 1.013122918778594
$ws.Cells.Item(7, 6).Value = 1.01662814436822
$ws.Cells.Item(7, 9).Value = 1.023278591626716
$ws.Cells.Item(7, 10).Value = 1.015561327166179
$ws.Cells.Item(7, 11).Value = 1.015948655736438
$ws.Cells.Item(7, 12).Value = 1.015665439685322
$ws.Cells.Item(7, 13).Value = 1.019161400592074
$ws.Cells.Item(7, 14).Value = 1.00935262377123

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.009113109594213
$ws.Cells.Item(8, 4).Value = 1.011728568650955
$ws.Cells.Item(8, 5).Value = 1.011590204848559
$ws.Cells.Item(8, 6).Value = 1.013204563420708
$ws.Cells.Item(8, 9).Value = 1.022925408098077
$ws.Cells.Item(8, 10).Value = 1.014252660193518
$ws.Cells.Item(8, 11).Value = 1.014529679201516
$ws.Cells.Item(8, 12).Value = 1.014391722538553
$ws.Cells.Item(8, 13).Value = 1.016001338343124
$ws.Cells.Item(8, 14).Value = 1.008922419093536

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.005881380305241
$ws.Cells.Item(9, 4).Value = 1.008705511637789
$ws.Cells.Item(9, 5).Value = 1.008833393751115
$ws.Cells.Item(9, 6).Value = 1.00678568889827
$ws.Cells.Item(9, 9).Value = 1.022200423623341
$ws.Cells.Item(9, 10).Value = 1.011875878842454
$ws.Cells.Item(9, 11).Value = 1.011959776244926
$ws.Cells.Item(9, 12).Value = 1.012087216056089
$ws.Cells.Item(9, 13).Value = 1.010046608407749
$ws.Cells.Item(9, 14).Value = 1.008140974237771

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.003678250791602
$ws.Cells.Item(10, 4).Value = 1.006643541996946
$ws.Cells.Item(10, 5).Value = 1.006955743419767
$ws.Cells.Item(10, 6).Value = 1.002234127190617
$ws.Cells.Item(10, 9).Value = 1.021646160605932
$ws.Cells.Item(10, 10).Value = 1.010241460384518
$ws.Cells.Item(10, 11).Value = 1.010197414604496
$ws.Cells.Item(10, 12).Value = 1.010508427753418
$ws.Cells.Item(10, 13).Value = 1.005804870230878
$ws.Cells.Item(10, 14).Value = 1.007603531192753

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.002712108930297
$ws.Cells.Item(11, 4).Value = 1.005739032984408
$ws.Cells.Item(11, 5).Value = 1.006132742772642
$ws.Cells.Item(11, 6).Value = 1.000195541274163
$ws.Cells.Item(11, 9).Value = 1.021388843948345
$ws.Cells.Item(11, 10).Value = 1.009521332991206
$ws.Cells.Item(11, 11).Value = 1.009422067674186
$ws.Cells.Item(11, 12).Value = 1.009814219386273
$ws.Cells.Item(11, 13).Value = 1.003900658822674
$ws.Cells.Item(11, 14).Value = 1.00736671590633

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.002351354576249
$ws.Cells.Item(12, 4).Value = 1.005401250502746
$ws.Cells.Item(12, 5).Value = 1.005825498246155
$ws.Cells.Item(12, 6).Value = 0.9994278554538939
$ws.Cells.Item(12, 9).Value = 1.021290617304148
$ws.Cells.Item(12, 10).Value = 1.009251927876413
$ws.Cells.Item(12, 11).Value = 1.009132177865306
$ws.Cells.Item(12, 12).Value = 1.009554722779117
$ws.Cells.Item(12, 13).Value = 1.003182929600929
$ws.Cells.Item(12, 14).Value = 1.007278119015379

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.002428823930299
$ws.Cells.Item(13, 4).Value = 1.005473788716769
$ws.Cells.Item(13, 5).Value = 1.005891473983286
$ws.Cells.Item(13, 6).Value = 0.9995930050484571
$ws.Cells.Item(13, 9).Value = 1.021311807850252
$ws.Cells.Item(13, 10).Value = 1.009309803855467
$ws.Cells.Item(13, 11).Value = 1.009194446704614
$ws.Cells.Item(13, 12).Value = 1.009610460532091
$ws.Cells.Item(13, 13).Value = 1.003337361165183
$ws.Cells.Item(13, 14).Value = 1.007297152291559

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.002682327687129
$ws.Cells.Item(14, 4).Value = 1.005711148985741
$ws.Cells.Item(14, 5).Value = 1.006107377660579
$ws.Cells.Item(14, 6).Value = 1.000132299330434
$ws.Cells.Item(14, 9).Value = 1.02138077879649
$ws.Cells.Item(14, 10).Value = 1.009499103319021
$ws.Cells.Item(14, 11).Value = 1.009398144199291
$ws.Cells.Item(14, 12).Value = 1.009792802985549
$ws.Cells.Item(14, 13).Value = 1.00384154534539
$ws.Cells.Item(14, 14).Value = 1.007359405480006

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.002838267990452
$ws.Cells.Item(15, 4).Value = 1.005857153244523
$ws.Cells.Item(15, 5).Value = 1.006240196760964
$ws.Cells.Item(15, 6).Value = 1.000463180757737
$ws.Cells.Item(15, 9).Value = 1.021422921776436
$ws.Cells.Item(15, 10).Value = 1.00961548121031
$ws.Cells.Item(15, 11).Value = 1.009523396651916
$ws.Cells.Item(15, 12).Value = 1.0099049318578
$ws.Cells.Item(15, 13).Value = 1.004150800438891
$ws.Cells.Item(15, 14).Value = 1.00739767728778

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.003742109163916
$ws.Cells.Item(16, 4).Value = 1.006703320917665
$ws.Cells.Item(16, 5).Value = 1.007010149252668
$ws.Cells.Item(16, 6).Value = 1.002367970247062
$ws.Cells.Item(16, 9).Value = 1.021662869162574
$ws.Cells.Item(16, 10).Value = 1.010288986907132
$ws.Cells.Item(16, 11).Value = 1.010248609601503
$ws.Cells.Item(16, 12).Value = 1.010554273188587
$ws.Cells.Item(16, 13).Value = 1.005929800836464
$ws.Cells.Item(16, 14).Value = 1.007619160030419

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.004305766535598
$ws.Cells.Item(17, 4).Value = 1.007230939158178
$ws.Cells.Item(17, 5).Value = 1.007490418954235
$ws.Cells.Item(17, 6).Value = 1.003544463576809
$ws.Cells.Item(17, 9).Value = 1.021808714138521
$ws.Cells.Item(17, 10).Value = 1.01070809965303
$ws.Cells.Item(17, 11).Value = 1.01070020473215
$ws.Cells.Item(17, 12).Value = 1.010958722507577
$ws.Cells.Item(17, 13).Value = 1.007027456233589
$ws.Cells.Item(17, 14).Value = 1.007756981028552

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.004633367978292
$ws.Cells.Item(18, 4).Value = 1.007537568617639
$ws.Cells.Item(18, 5).Value = 1.007769594027454
$ws.Cells.Item(18, 6).Value = 1.00422417784475
$ws.Cells.Item(18, 9).Value = 1.021892115505348
$ws.Cells.Item(18, 10).Value = 1.010951367063032
$ws.Cells.Item(18, 11).Value = 1.010962436062044
$ws.Cells.Item(18, 12).Value = 1.011193613208824
$ws.Cells.Item(18, 13).Value = 1.007661206017583
$ws.Cells.Item(18, 14).Value = 1.007836975432551

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.004744874776143
$ws.Cells.Item(19, 4).Value = 1.007641932764769
$ws.Cells.Item(19, 5).Value = 1.007864624412054
$ws.Cells.Item(19, 6).Value = 1.004454846334689
$ws.Cells.Item(19, 9).Value = 1.021920271712059
$ws.Cells.Item(19, 10).Value = 1.011034113991051
$ws.Cells.Item(19, 11).Value = 1.011051652189238
$ws.Cells.Item(19, 12).Value = 1.011273533611196
$ws.Cells.Item(19, 13).Value = 1.007876204930487
$ws.Cells.Item(19, 14).Value = 1.00786418509266

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.004245412894813
$ws.Cells.Item(20, 4).Value = 1.007174447088361
$ws.Cells.Item(20, 5).Value = 1.007438990001944
$ws.Cells.Item(20, 6).Value = 1.00341891281462
$ws.Cells.Item(20, 9).Value = 1.021793239156063
$ws.Cells.Item(20, 10).Value = 1.010663256689616
$ws.Cells.Item(20, 11).Value = 1.010651874912773
$ws.Cells.Item(20, 12).Value = 1.010915434521996
$ws.Cells.Item(20, 13).Value = 1.006910361852493
$ws.Cells.Item(20, 14).Value = 1.007742235040907

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.002607729699117
$ws.Cells.Item(21, 4).Value = 1.005641302663083
$ws.Cells.Item(21, 5).Value = 1.006043842473613
$ws.Cells.Item(21, 6).Value = 0.9999737818872544
$ws.Cells.Item(21, 9).Value = 1.021360542051932
$ws.Cells.Item(21, 10).Value = 1.009443412725443
$ws.Cells.Item(21, 11).Value = 1.009338213025364
$ws.Cells.Item(21, 12).Value = 1.00973915325026
$ws.Cells.Item(21, 13).Value = 1.003693365648453
$ws.Cells.Item(21, 14).Value = 1.007341091090635

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.001567108860067
$ws.Cells.Item(22, 4).Value = 1.004666865853775
$ws.Cells.Item(22, 5).Value = 1.005157690225753
$ws.Cells.Item(22, 6).Value = 0.9977469910722081
$ws.Cells.Item(22, 9).Value = 1.021073143483512
$ws.Cells.Item(22, 10).Value = 1.008665324280023
$ws.Cells.Item(22, 11).Value = 1.008501289250969
$ws.Cells.Item(22, 12).Value = 1.008990083600116
$ws.Cells.Item(22, 13).Value = 1.001610272144382
$ws.Cells.Item(22, 14).Value = 1.007085203336629

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.002119819223535
$ws.Cells.Item(23, 4).Value = 1.005184446832403
$ws.Cells.Item(23, 5).Value = 1.005628323248888
$ws.Cells.Item(23, 6).Value = 0.9989333102599026
$ws.Cells.Item(23, 9).Value = 1.021226969748428
$ws.Cells.Item(23, 10).Value = 1.009078876663895
$ws.Cells.Item(23, 11).Value = 1.008946017266213
$ws.Cells.Item(23, 12).Value = 1.009388096038074
$ws.Cells.Item(23, 13).Value = 1.002720386211757
$ws.Cells.Item(23, 14).Value = 1.007221208489272

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.004272687720373
$ws.Cells.Item(24, 4).Value = 1.007199976884807
$ws.Cells.Item(24, 5).Value = 1.007462231489618
$ws.Cells.Item(24, 6).Value = 1.003475663910708
$ws.Cells.Item(24, 9).Value = 1.021800236784383
$ws.Cells.Item(24, 10).Value = 1.010683522980798
$ws.Cells.Item(24, 11).Value = 1.010673716709996
$ws.Cells.Item(24, 12).Value = 1.010934997641495
$ws.Cells.Item(24, 13).Value = 1.006963291811816
$ws.Cells.Item(24, 14).Value = 1.007748899335216

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.006725234458547
$ws.Cells.Item(25, 4).Value = 1.009495066429831
$ws.Cells.Item(25, 5).Value = 1.009552939992108
$ws.Cells.Item(25, 6).Value = 1.00849212760781
$ws.Cells.Item(25, 9).Value = 1.02240020001346
$ws.Cells.Item(25, 10).Value = 1.012498947110091
$ws.Cells.Item(25, 11).Value = 1.012632626159716
$ws.Cells.Item(25, 12).Value = 1.012690307311561
$ws.Cells.Item(25, 13).Value = 1.011633025626726
$ws.Cells.Item(25, 14).Value = 1.008345841479275
